# Update NATMI ligand-receptor statistics for Vegfb-Flt1 sheet
# following revised ligand/receptor-expressing cell counts (1 -> 3),
# per author's note: "Natmi following Dr Hou advice".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.625493333333333
$ws.Range("H2").Value = 4.87648
$ws.Range("I2").Value = 0.1468796758507528
$ws.Range("J2").Value = 0.1468796758507528
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 61.629167
$ws.Range("N2").Value = 184.887501
$ws.Range("O2").Value = 0.8452417044501688
$ws.Range("P2").Value = 0.8452417044501688
$ws.Range("Q2").Value = 100.1778000973867
$ws.Range("R2").Value = 901.60020087648
$ws.Range("S2").Value = 0.1241488275651786
$ws.Range("T2").Value = 0.1241488275651786

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.625493333333333
$ws.Range("H3").Value = 4.87648
$ws.Range("I3").Value = 0.1468796758507528
$ws.Range("J3").Value = 0.1468796758507528
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.357683666666667
$ws.Range("N3").Value = 4.073051
$ws.Range("O3").Value = 0.01862058035795759
$ws.Range("P3").Value = 0.01862058035795759
$ws.Range("Q3").Value = 2.206905748942222
$ws.Range("R3").Value = 19.86215174048
$ws.Range("S3").Value = 0.002734984807129704
$ws.Range("T3").Value = 0.002734984807129705

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.625493333333333
$ws.Range("H4").Value = 4.87648
$ws.Range("I4").Value = 0.1468796758507528
$ws.Range("J4").Value = 0.1468796758507528
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.926218666666665
$ws.Range("N4").Value = 29.778656
$ws.Range("O4").Value = 0.1361377151918736
$ws.Range("P4").Value = 0.1361377151918736
$ws.Range("Q4").Value = 16.13500226787555
$ws.Range("R4").Value = 145.21502041088
$ws.Range("S4").Value = 0.01999586347844449
$ws.Range("T4").Value = 0.0199958634784445

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.041193333333333
$ws.Range("H5").Value = 15.12358
$ws.Range("I5").Value = 0.4555225343081337
$ws.Range("J5").Value = 0.4555225343081337
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 61.629167
$ws.Range("N5").Value = 184.887501
$ws.Range("O5").Value = 0.8452417044501688
$ws.Range("P5").Value = 0.8452417044501688
$ws.Range("Q5").Value = 310.6845458192867
$ws.Range("R5").Value = 2796.16091237358
$ws.Range("S5").Value = 0.3850266433140674
$ws.Range("T5").Value = 0.3850266433140674

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.041193333333333
$ws.Range("H6").Value = 15.12358
$ws.Range("I6").Value = 0.4555225343081337
$ws.Range("J6").Value = 0.4555225343081337
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.357683666666667
$ws.Range("N6").Value = 4.073051
$ws.Range("O6").Value = 0.01862058035795759
$ws.Range("P6").Value = 0.01862058035795759
$ws.Range("Q6").Value = 6.844345849175554
$ws.Range("R6").Value = 61.59911264258
$ws.Range("S6").Value = 0.008482093954945095
$ws.Range("T6").Value = 0.008482093954945097

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.041193333333333
$ws.Range("H7").Value = 15.12358
$ws.Range("I7").Value = 0.4555225343081337
$ws.Range("J7").Value = 0.4555225343081337
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.926218666666665
$ws.Range("N7").Value = 29.778656
$ws.Range("O7").Value = 0.1361377151918736
$ws.Range("P7").Value = 0.1361377151918736
$ws.Range("Q7").Value = 50.03998736760888
$ws.Range("R7").Value = 450.35988630848
$ws.Range("S7").Value = 0.06201379703912115
$ws.Range("T7").Value = 0.06201379703912117

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.400149666666667
$ws.Range("H8").Value = 13.200449
$ws.Range("I8").Value = 0.3975977898411136
$ws.Range("J8").Value = 0.3975977898411136
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 61.629167
$ws.Range("N8").Value = 184.887501
$ws.Range("O8").Value = 0.8452417044501688
$ws.Range("P8").Value = 0.8452417044501688
$ws.Range("Q8").Value = 271.1775586319943
$ws.Range("R8").Value = 2440.598027687949
$ws.Range("S8").Value = 0.3360662335709229
$ws.Range("T8").Value = 0.3360662335709229

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.400149666666667
$ws.Range("H9").Value = 13.200449
$ws.Range("I9").Value = 0.3975977898411136
$ws.Range("J9").Value = 0.3975977898411136
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.357683666666667
$ws.Range("N9").Value = 4.073051
$ws.Range("O9").Value = 0.01862058035795759
$ws.Range("P9").Value = 0.01862058035795759
$ws.Range("Q9").Value = 5.97401133332211
$ws.Range("R9").Value = 53.76610199989899
$ws.Range("S9").Value = 0.007403501595882789
$ws.Range("T9").Value = 0.007403501595882789

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.400149666666667
$ws.Range("H10").Value = 13.200449
$ws.Range("I10").Value = 0.3975977898411136
$ws.Range("J10").Value = 0.3975977898411136
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.926218666666665
$ws.Range("N10").Value = 29.778656
$ws.Range("O10").Value = 0.1361377151918736
$ws.Range("P10").Value = 0.1361377151918736
$ws.Range("Q10").Value = 43.67684775739377
$ws.Range("R10").Value = 393.091629816544
$ws.Range("S10").Value = 0.05412805467430792
$ws.Range("T10").Value = 0.05412805467430794
